$wb = $excel.ActiveWorkbook

$updates = @{
    "F4"  = 1005
    "F7"  = 2532
    "F9"  = 1574
    "F11" = 180
    "F13" = 483
    "F15" = 40
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
